$d = $word.ActiveDocument

# The paragraph ending "...already had this sort of constructor defined." is
# the anchor: the three new explanatory paragraphs about MarieSpawnBehaviour
# are inserted right after it, before the _GoBack bookmark paragraph.
$anchor = $d.Paragraphs(4)
$anchorRange = $anchor.Range

# Insert three empty paragraphs first (in document order) so that none of
# them inadvertently inherit the tab stop we add to the first one afterwards.
$anchorRange.InsertParagraphAfter()
$p1 = $d.Paragraphs(5)

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(6)

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(7)

# --- New paragraph 1: "The class MamboMarie..." ---
$p1.Range.Text = 'The class MamboMarie, like other classes such as Zombie and Farmer uses an array to hold a set of behaviours. These behaviours are MarieSpawnBehaviour, ChantBehaviour, and WanderBehaviour.'

# --- New paragraph 2: "In order to implement MamboMarie..." ---
$p2.Range.Text = 'In order to implement MamboMarie appearing 5 percent of the time, the behaviour MarieSpawnBehaviour is created. This class implements the Behaviour interface. By implementing an interface, we can easily attain methods that we will need in order to develop this behaviour such as the getAction() method. This is good practice as this interface requires that all methods that implement it are required to have it’s methods, so we can uphold  consistency amongst all Behaviours.'

# --- New paragraph 3: "The MarieSpawnBehaviour class implements..." ---
$p3.Range.Text = 'The MarieSpawnBehaviour class implements the getAction method which returns a MarieSpawnAction, or null if it is not. A helper class to determine a random edge location getRandEdgeLocation is utilized which determines and returns a random Location on the top edge of the map. This method is called in the getAction action method. An alternative for this would be to determine the random location inside of getAction, however this would be bad design practice as this means that the code is not modularised and therefore more difficult to interpret. In this manner, using a helper method makes the code more readable.'

# Only the first new paragraph carries the left tab stop at 6824 twips
# (341.2pt, since TabStops.Add takes points).
$p1.TabStops.Add(341.2)

# Remove one of the two blank paragraphs that followed the _GoBack bookmark
# paragraph (there were two `<w:p/>` in a row; only one should remain).
$blank1 = $d.Paragraphs(9)
$blank1.Range.Delete()

Write-Output "Done"
